# Update automàtic: dades i banners [2026-02-20 11:38]
$wb = $excel.ActiveWorkbook

$wsDades = $wb.Worksheets.Item("Dades_Període")
$wsCapcaleres = $wb.Worksheets.Item("Estudi_Capçaleres")

# Row 2 (station XJ, current period 10:30 - 11:00)
$wsDades.Range("H2").Value = "2026-02-20 11:38:39"
$wsDades.Range("I2").Value = "11:00"
$wsDades.Range("J2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T11:00Z"

# Rows 3-6 (other period rows) - only the extraction timestamp changes
$wsDades.Range("H3").Value = "2026-02-20 11:38:41"
$wsDades.Range("H4").Value = "2026-02-20 11:38:41"
$wsDades.Range("H5").Value = "2026-02-20 11:38:41"
$wsDades.Range("H6").Value = "2026-02-20 11:38:41"

# Estudi_Capçaleres sheet: banner URL updated to new period
$wsCapcaleres.Range("F2").Value = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T11:00Z"
